$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Myles feedback on 10 Pager v1
$ws.Range("A3").Value = "10 Pager v1"
$ws.Range("B3").Value = "Myles"
$ws.Range("C3").Value = "•With the Target Market, I'd like to see at least 5 profiles from social media where connections have been made and characteristics identified. You have identified the characteristics but I don't know where you got the info from.
•I'd also like to see wire-frame mockups for the UI and mechanics."

# Row 4 - Elric feedback on 10 Pager v1
$ws.Range("A4").Value = "10 Pager v1"
$ws.Range("B4").Value = "Elric"

# Update existing feedback response (row 2, column D)
$ws.Range("D2").Value = "•Flesh out the theme and message I'm trying to tell with the game, and fit in the gameplay."

$ws.Range("C4").Value = "•What platform is it on? (list in overview somewhere).
•Maybe the japanese culture could be more of a flourish or emphasis rather than a core part of the aesthetic, so as not to confuse the player? OR make that blend of cultures a lot clearer in the story.
•Images worked really well to explain the look and environmental mechanics!
•What is the age rating.
•Clarify if each combat mechanic is consistent for NPCs and PC enemy or if some of them are only intended for one or the other."
$ws.Range("D4").Value = "•Add information such as platform and age rating.
•Explain exactly how the culture blend will be used and what it adds to the game.
•Add more info for mechanics"

$ws.Range("D3").Value = "•Add profiles for target market and description of relationship.
•Finish wireframe."

# Apply styles/formats matching row 2 (vertical top, wrap text where applicable)
$ws.Range("A3:B4").VerticalAlignment = -4160
$ws.Range("C3:C4").WrapText = $true
$ws.Range("D3:D4").VerticalAlignment = -4160
$ws.Range("D3:D4").WrapText = $true

# Row heights to match the wrapped-text autofit result in the source workbook
$ws.Rows.Item(3).RowHeight = 105
$ws.Rows.Item(4).RowHeight = 165

# Update selection to match diff's sheetView selection
$ws.Range("E11").Select()
